# Working preview in dashboard:
# Remove the stale/placeholder rows from the shelter data sheet and
# refresh the "Active" / "Status" flags for the rows that remain so the
# dashboard preview reflects the current facility survey.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that no longer belong in the sheet: "wdw", "F. Mendoza
# Memorial Elem Sch.", "Calumpit Sports Complex", "Gatbuca Basketball
# Court", "Dona Damiana Elem School", "Meysulao Multipurpose/E.C." and
# "Calizon Dike" (original rows 2,3,4,5,7,9,10). Delete from the bottom
# up so earlier row numbers stay valid as later ones are removed.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

# After the deletions, refresh Active/Status flags on the surviving rows
# (now at rows 2-4): San Miguel Meysulao High School, Danga Dike and
# San Marcos Elem. Sch.
$ws.Range("A2").Value = $true
$ws.Range("A3").Value = $true
$ws.Range("L3").Value = "Partially Built"
$ws.Range("A4").Value = $true
